$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 83

# Columns A-D hold text that Excel would otherwise auto-convert (dates,
# leading-zero week numbers). Force text number format first, then clear
# the formatting afterwards so the new row ends up unstyled just like the
# rest of the data rows.
$textRange = $ws.Range("A$($row):D$($row)")
$textRange.NumberFormat = "@"

$ws.Range("A$($row)").Value = "2024-01-22"
$ws.Range("B$($row)").Value = "12:48:03"
$ws.Range("C$($row)").Value = "Monday"
$ws.Range("D$($row)").Value = "03"

$textRange.ClearFormats()

$ws.Range("E$($row)").Value = 138561
$ws.Range("F$($row)").Value = 141159
$ws.Range("G$($row)").Value = 171140
$ws.Range("H$($row)").Value = 148646
$ws.Range("I$($row)").Value = -1
$ws.Range("J$($row)").Value = 123068
$ws.Range("K$($row)").Value = 223726
$ws.Range("L$($row)").Value = 255734
$ws.Range("M$($row)").Value = 185596
$ws.Range("N$($row)").Value = 110408
$ws.Range("O$($row)").Value = 41309
$ws.Range("P$($row)").Value = 30897
$ws.Range("Q$($row)").Value = 73605
$ws.Range("R$($row)").Value = -1
$ws.Range("S$($row)").Value = 42600
$ws.Range("T$($row)").Value = -1
